$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1598.1852
$ws.Range("J28").Value = 6198.6
$ws.Range("L28").Value = 6198.6
$ws.Range("N28").Value = -7168.6
$ws.Range("H43").Value = 3715.4614
$ws.Range("I43").Value = 4889.4
$ws.Range("J43").Value = 2981.75
$ws.Range("K43").Value = 4889.4
$ws.Range("L43").Value = 2981.75
$ws.Range("M43").Value = -4820.4
$ws.Range("N43").Value = -3119.75
$ws.Range("H98").Value = 2637.8064
$ws.Range("I98").Value = 2693.5667
$ws.Range("K98").Value = 2693.5667
$ws.Range("M98").Value = -1195.5667
$ws.Range("H116").Value = 2344.3076
$ws.Range("I116").Value = 2343.182
$ws.Range("K116").Value = 2343.182
$ws.Range("M116").Value = 1098.818
$ws.Range("H122").Value = 2637.8064
$ws.Range("I122").Value = 2693.5667
$ws.Range("K122").Value = 8080.7001
$ws.Range("M122").Value = -5630.7001
$ws.Range("H132").Value = 8349.789000000001
$ws.Range("I132").Value = 9275.706
$ws.Range("J132").Value = 479.5
$ws.Range("K132").Value = 27827.118
$ws.Range("L132").Value = 1438.5
$ws.Range("M132").Value = -25297.118
$ws.Range("N132").Value = -6498.5
$ws.Range("H137").Value = 1530.0834
$ws.Range("I137").Value = 1557.7188
$ws.Range("K137").Value = 4673.1564
$ws.Range("M137").Value = -2123.1564
$ws.Range("H138").Value = 351938.2
$ws.Range("I138").Value = 5716.3335
$ws.Range("J138").Value = 443049.22
$ws.Range("K138").Value = 17149.0005
$ws.Range("L138").Value = 1329147.66
$ws.Range("M138").Value = -12009.0005
$ws.Range("N138").Value = -1339427.66
$ws.Range("H141").Value = 818.3333
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 32929.145
$ws.Range("J45").Value = 4018
$ws.Range("L45").Value = 4018
$ws.Range("N45").Value = -4772
$ws.Range("H48").Value = 299999
$ws.Range("J48").Value = 299999
$ws.Range("L48").Value = 299999
$ws.Range("N48").Value = -300767
$ws.Range("H61").Value = 3229.2068
$ws.Range("I61").Value = 2052.2273
$ws.Range("K61").Value = 2052.2273
$ws.Range("M61").Value = -1840.2273
$ws.Range("H74").Value = 234516.05
$ws.Range("I74").Value = 280593.6
$ws.Range("K74").Value = 280593.6
$ws.Range("M74").Value = -279719.6
$ws.Range("H77").Value = 234516.05
$ws.Range("I77").Value = 280593.6
$ws.Range("K77").Value = 1402968
$ws.Range("M77").Value = -1398600
$ws.Range("H136").Value = 3229.2068
$ws.Range("I136").Value = 2052.2273
$ws.Range("K136").Value = 6156.6819
$ws.Range("M136").Value = -3606.6819

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 5654557
$ws.Range("I105").Value = 296303.75
$ws.Range("K105").Value = 296303.75
$ws.Range("M105").Value = -294556.75

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3955.5938
$ws.Range("J31").Value = 5089.36
$ws.Range("L31").Value = 5089.36
$ws.Range("N31").Value = -5679.36
$ws.Range("H34").Value = 3955.5938
$ws.Range("J34").Value = 5089.36
$ws.Range("L34").Value = 5089.36
$ws.Range("N34").Value = -5493.36
$ws.Range("H86").Value = 4998.5
$ws.Range("J86").Value = 4998.5
$ws.Range("L86").Value = 4998.5
$ws.Range("N86").Value = -7244.5
$ws.Range("H89").Value = 4998.5
$ws.Range("J89").Value = 4998.5
$ws.Range("L89").Value = 24992.5
$ws.Range("N89").Value = -36224.5
$ws.Range("H134").Value = 2466.5833
$ws.Range("I134").Value = 2087
$ws.Range("J134").Value = 3605.3333
$ws.Range("K134").Value = 6261
$ws.Range("L134").Value = 10815.9999
$ws.Range("M134").Value = -3726
$ws.Range("N134").Value = -15885.9999
$ws.Range("H141").Value = 333285
$ws.Range("J141").Value = 333285
$ws.Range("L141").Value = 333285
$ws.Range("N141").Value = -343645

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 850.0714
$ws.Range("I2").Value = 59.05263
$ws.Range("J2").Value = 2520
$ws.Range("K2").Value = 354.31578
$ws.Range("L2").Value = 15120
$ws.Range("M2").Value = -241.31578
$ws.Range("N2").Value = -15346
$ws.Range("H37").Value = 18911696
$ws.Range("J37").Value = 18911696
$ws.Range("L37").Value = 56735088
$ws.Range("N37").Value = -56735312
$ws.Range("H107").Value = 899.4
$ws.Range("J107").Value = 961.0769
$ws.Range("L107").Value = 2883.2307
$ws.Range("N107").Value = -6723.2307
$ws.Range("H113").Value = 1622.5
$ws.Range("I113").Value = 871.625
$ws.Range("K113").Value = 2614.875
$ws.Range("M113").Value = -444.875
$ws.Range("H115").Value = 219810.5
$ws.Range("J115").Value = 256383.08
$ws.Range("L115").Value = 769149.24
$ws.Range("N115").Value = -771499.24
$ws.Range("H120").Value = 27133.334
$ws.Range("J120").Value = 30000
$ws.Range("L120").Value = 90000
$ws.Range("N120").Value = -99676
$ws.Range("H122").Value = 1557.6111
$ws.Range("J122").Value = 1535.8667
$ws.Range("L122").Value = 13822.8003
$ws.Range("N122").Value = -18722.8003

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2465.6
$ws.Range("I102").Value = 2157
$ws.Range("J102").Value = 3082.8
$ws.Range("K102").Value = 2157
$ws.Range("L102").Value = 3082.8
$ws.Range("M102").Value = -535
$ws.Range("N102").Value = -6326.8
$ws.Range("H132").Value = 3578
$ws.Range("I132").Value = 2945.4333
$ws.Range("J132").Value = 4632.278
$ws.Range("K132").Value = 8836.2999
$ws.Range("L132").Value = 13896.834
$ws.Range("M132").Value = -6306.2999
$ws.Range("N132").Value = -18956.834

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").ClearContents()
$ws.Range("N8").Value = 0
$ws.Range("H46").Value = 1640.1428
$ws.Range("I46").Value = 2308.3333
$ws.Range("J46").Value = 1501.8966
$ws.Range("K46").Value = 2308.3333
$ws.Range("L46").Value = 1501.8966
$ws.Range("M46").Value = -2120.3333
$ws.Range("N46").Value = -1877.8966
$ws.Range("H55").Value = 492.64285
$ws.Range("I55").Value = 465.33334
$ws.Range("K55").Value = 465.33334
$ws.Range("M55").Value = -292.33334
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("H119").Value = 49999
$ws.Range("J119").Value = 49999
$ws.Range("L119").Value = 49999
$ws.Range("N119").Value = -59675
$ws.Range("H132").Value = 4078.0417
$ws.Range("I132").Value = 3171.3635
$ws.Range("K132").Value = 9514.0905
$ws.Range("M132").Value = -6984.0905

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H122").Value = 11366410
$ws.Range("I122").Value = 2394.353
$ws.Range("K122").Value = 7183.059
$ws.Range("M122").Value = -4733.059
$ws.Range("I126").Value = 1913.75
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 5741.25
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -3271.25
$ws.Range("N126").Value = -8540
$ws.Range("H132").Value = 4746.278
$ws.Range("I132").Value = 4577.3
$ws.Range("J132").Value = 5591.1665
$ws.Range("K132").Value = 13731.9
$ws.Range("L132").Value = 16773.4995
$ws.Range("M132").Value = -11201.9
$ws.Range("N132").Value = -21833.4995

